# Added better date format sensing.
#
# - C1 header ("Notes") gets centered (like the other header cells already are).
# - B3's existing custom date format (numFmtId 166) gets right-aligned.
# - A new date, E3, is added using a brand-new custom format "ddd mm/d/yy".
# - C9 ("test123") gets centered within its shaded cell.
# - Selection cursor moves to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center the "Notes" header cell (C1)
$ws.Range("C1").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# Right-align the already-custom-formatted date in B3
$ws.Range("B3").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignRight

# New date value in E3 with a new custom number format (ddd mm/d/yy)
$ws.Range("E3").Value = 42502
$ws.Range("E3").NumberFormat = 'ddd\ mm/d/yy'

# Center the text in the shaded C9 cell
$ws.Range("C9").HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# Move the active selection to C3
[void]$ws.Range("C3").Select()
